$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Rename the "ID" column header to "Id" (cascades into the table definition
# and the shared-strings table exactly like the authored change).
$ws.Range("A1").Value = "Id"

# Move the active selection back to A2 (previously left on U17).
$ws.Range("A2").Select()
